$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# row 12
$ws.Range("H12").Value = 272.18182
$ws.Range("I12").Value = 287.125
$ws.Range("K12").Value = 287.125
$ws.Range("M12").Value = -117.125
# row 40
$ws.Range("H40").Value = 3578.7856
$ws.Range("I40").Value = 2343.2856
$ws.Range("J40").Value = 4814.2856
$ws.Range("K40").Value = 2343.2856
$ws.Range("L40").Value = 4814.2856
$ws.Range("M40").Value = -2168.2856
$ws.Range("N40").Value = -5164.2856
# row 112
$ws.Range("H112").Value = 1358.6
$ws.Range("J112").Value = 1398.375
$ws.Range("L112").Value = 4195.125
$ws.Range("N112").Value = -6411.125
# row 116
$ws.Range("H116").Value = 15000
$ws.Range("I116").Value = 14500
$ws.Range("J116").Value = 15500
$ws.Range("K116").Value = 14500
$ws.Range("L116").Value = 15500
$ws.Range("M116").Value = -11058
$ws.Range("N116").Value = -22384
# row 125
$ws.Range("H125").Value = 70666.336
$ws.Range("I125").Value = 0
$ws.Range("J125").Value = 70666.336
$ws.Range("K125").Value = 0
$ws.Range("L125").ClearContents()
$ws.Range("M125").Value = 635997.024
$ws.Range("N125").Value = -640917.024
# row 132
$ws.Range("H132").Value = 6676629
$ws.Range("I132").Value = 3802.5386
$ws.Range("J132").Value = 50050000
$ws.Range("K132").Value = 11407.6158
$ws.Range("L132").Value = 150150000
$ws.Range("M132").Value = -8877.6158
$ws.Range("N132").Value = -150155060
# row 137
$ws.Range("H137").Value = 4348947
$ws.Range("I137").Value = 1451.3334
$ws.Range("J137").Value = 10870190
$ws.Range("K137").Value = 4354.0002
$ws.Range("L137").Value = 32610570
$ws.Range("M137").Value = -1804.0002
$ws.Range("N137").Value = -32615670

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# row 74
$ws.Range("H74").Value = 3269.9412
$ws.Range("I74").Value = 2965.6
$ws.Range("J74").Value = 3704.7144
$ws.Range("K74").Value = 2965.6
$ws.Range("L74").Value = 3704.7144
$ws.Range("M74").Value = -2091.6
$ws.Range("N74").Value = -5452.7144
# row 77
$ws.Range("H77").Value = 3269.9412
$ws.Range("I77").Value = 2965.6
$ws.Range("J77").Value = 3704.7144
$ws.Range("K77").Value = 14828
$ws.Range("L77").Value = 18523.572
$ws.Range("M77").Value = -10460
$ws.Range("N77").Value = -27259.572
# row 132
$ws.Range("H132").Value = 7431.55
$ws.Range("I132").Value = 7719
$ws.Range("J132").Value = 1970
$ws.Range("K132").Value = 23157
$ws.Range("L132").Value = 5910
$ws.Range("M132").Value = -20627
$ws.Range("N132").Value = -10970

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# row 55
$ws.Range("H55").Value = 45832.918
$ws.Range("J55").Value = 45832.918
$ws.Range("L55").Value = 45832.918
$ws.Range("N55").Value = -46378.918
# row 76
$ws.Range("H76").Value = 24322
$ws.Range("J76").Value = 24322
$ws.Range("L76").Value = 24322
$ws.Range("N76").Value = -24952
# row 79
$ws.Range("H79").Value = 24322
$ws.Range("J79").Value = 24322
$ws.Range("L79").Value = 24322
$ws.Range("N79").Value = -26506
# row 86
$ws.Range("H86").Value = 35723024
$ws.Range("I86").Value = 13848.875
$ws.Range("J86").Value = 83335256
$ws.Range("K86").Value = 13848.875
$ws.Range("L86").Value = 83335256
$ws.Range("M86").Value = -12725.875
$ws.Range("N86").Value = -83337502
# row 89
$ws.Range("H89").Value = 35723024
$ws.Range("I89").Value = 13848.875
$ws.Range("J89").Value = 83335256
$ws.Range("K89").Value = 69244.375
$ws.Range("L89").Value = 416676280
$ws.Range("M89").Value = -63628.375
$ws.Range("N89").Value = -416687512
# row 105
$ws.Range("H105").Value = 13756.556
$ws.Range("I105").Value = 18468.166
$ws.Range("K105").Value = 18468.166
$ws.Range("M105").Value = -16721.166

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# row 58
$ws.Range("H58").Value = 2314.9644
$ws.Range("I58").Value = 2268.0386
$ws.Range("K58").Value = 2268.0386
$ws.Range("M58").Value = -2065.0386
# row 96
$ws.Range("H96").Value = 19372.857
$ws.Range("J96").Value = 20101.834
$ws.Range("L96").Value = 20101.834
$ws.Range("N96").Value = -25593.834
# row 97
$ws.Range("H97").Value = 115965
$ws.Range("J97").Value = 115965
$ws.Range("L97").Value = 115965
$ws.Range("N97").Value = -117947
# row 132
$ws.Range("H132").Value = 3806.2307
$ws.Range("I132").Value = 2589.182
$ws.Range("K132").Value = 7767.545999999999
$ws.Range("M132").Value = -5237.545999999999
# row 136
$ws.Range("H136").Value = 2314.9644
$ws.Range("I136").Value = 2268.0386
$ws.Range("K136").Value = 6804.1158
$ws.Range("M136").Value = -4254.1158

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# row 129
$ws.Range("H129").Value = 4693.4546
$ws.Range("J129").Value = 9100.25
$ws.Range("L129").Value = 27300.75
$ws.Range("N129").Value = -37300.75
# row 134
$ws.Range("H134").Value = 4761
$ws.Range("I134").Value = 1797.8334
$ws.Range("K134").Value = 5393.5002
$ws.Range("M134").Value = -323.5002000000004

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# row 126
$ws.Range("H126").Value = 3059
$ws.Range("I126").Value = 1777.909
$ws.Range("J126").Value = 5877.4
$ws.Range("K126").Value = 5333.727000000001
$ws.Range("L126").Value = 17632.2
$ws.Range("M126").Value = -2863.727000000001
$ws.Range("N126").Value = -22572.2

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# row 2
$ws.Range("H2").Value = 392.96667
$ws.Range("I2").Value = 392.96667
$ws.Range("K2").Value = 392.96667
$ws.Range("M2").Value = -280.96667
# row 22
$ws.Range("H22").Value = 3708.0908
$ws.Range("I22").Value = 799
$ws.Range("K22").Value = 799
$ws.Range("M22").Value = -504
# row 27
$ws.Range("H27").Value = 3708.0908
$ws.Range("I27").Value = 799
$ws.Range("K27").Value = 799
$ws.Range("M27").Value = -692
# row 46
$ws.Range("H46").Value = 1255.75
$ws.Range("J46").Value = 2577.111
$ws.Range("L46").Value = 2577.111
$ws.Range("N46").Value = -2953.111
# row 55
$ws.Range("H55").Value = 3027.75
$ws.Range("I55").Value = 2556.5
$ws.Range("J55").Value = 3499
$ws.Range("K55").Value = 2556.5
$ws.Range("L55").Value = 3499
$ws.Range("M55").Value = -2383.5
$ws.Range("N55").Value = -3845
# row 61
$ws.Range("H61").Value = 5292.4736
$ws.Range("J61").Value = 9777.666999999999
$ws.Range("L61").Value = 9777.666999999999
$ws.Range("N61").Value = -10181.667
# row 74
$ws.Range("H74").Value = 500023520
$ws.Range("I74").Value = 40598.5
$ws.Range("K74").Value = 40598.5
$ws.Range("M74").Value = -39600.5
# row 77
$ws.Range("H77").Value = 500023520
$ws.Range("I77").Value = 40598.5
$ws.Range("K77").Value = 121795.5
$ws.Range("M77").Value = -116803.5
# row 93
$ws.Range("H93").Value = 4545.3335
$ws.Range("I93").Value = 2428.5715
$ws.Range("J93").Value = 6105.0527
$ws.Range("K93").Value = 2428.5715
$ws.Range("L93").Value = 6105.0527
$ws.Range("M93").Value = -1180.5715
$ws.Range("N93").Value = -8601.0527
# row 113
$ws.Range("H113").Value = 5292.4736
$ws.Range("J113").Value = 9777.666999999999
$ws.Range("L113").Value = 9777.666999999999
$ws.Range("N113").Value = -14117.667
# row 132
$ws.Range("H132").Value = 4914.3213
$ws.Range("I132").Value = 3708.9546
$ws.Range("J132").Value = 9334
$ws.Range("K132").Value = 11126.8638
$ws.Range("L132").Value = 28002
$ws.Range("M132").Value = -8596.863799999999
$ws.Range("N132").Value = -33062
# row 136
$ws.Range("H136").Value = 20853
$ws.Range("I136").Value = 3874
$ws.Range("J136").Value = 64072.273
$ws.Range("K136").Value = 11622
$ws.Range("L136").Value = 192216.819
$ws.Range("M136").Value = -9072
$ws.Range("N136").Value = -197316.819

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# row 2
$ws.Range("H2").Value = 2070815.2
$ws.Range("I2").Value = 3450200.5
$ws.Range("K2").Value = 3450200.5
$ws.Range("M2").Value = -3450088.5
# row 26
$ws.Range("H26").Value = 0
$ws.Range("I26").Value = 0
$ws.Range("K26").Value = 0
$ws.Range("M26").ClearContents()
# row 29
$ws.Range("H29").Value = 866.3333
$ws.Range("I29").Value = 500
$ws.Range("J29").Value = 1049.5
$ws.Range("K29").Value = 500
$ws.Range("L29").Value = 1049.5
$ws.Range("M29").Value = -210
$ws.Range("N29").Value = -1629.5
# row 47
$ws.Range("H47").Value = 11069
$ws.Range("J47").Value = 11069
$ws.Range("L47").Value = 11069
$ws.Range("N47").Value = -12213
# row 80
$ws.Range("H80").Value = 13949
$ws.Range("J80").Value = 13949
$ws.Range("L80").Value = 13949
$ws.Range("N80").Value = -15945
# row 83
$ws.Range("H83").Value = 13949
$ws.Range("J83").Value = 13949
$ws.Range("L83").Value = 41847
$ws.Range("N83").Value = -51831
# row 107
$ws.Range("H107").Value = 50000670
$ws.Range("J107").Value = 100000590
$ws.Range("L107").Value = 300001770
$ws.Range("N107").Value = -300005610
# row 122
$ws.Range("H122").Value = 607054.4
$ws.Range("I122").Value = 975011
$ws.Range("J122").Value = 9124.875
$ws.Range("K122").Value = 2925033
$ws.Range("L122").Value = 27374.625
$ws.Range("M122").Value = -2922583
$ws.Range("N122").Value = -32274.625
# row 126
$ws.Range("H126").Value = 9260617
$ws.Range("I126").Value = 1245.75
$ws.Range("K126").Value = 3737.25
$ws.Range("M126").Value = -1267.25
# row 132
$ws.Range("H132").Value = 12829736
$ws.Range("I132").Value = 12829736
$ws.Range("K132").Value = 38489208
$ws.Range("M132").Value = -38486678
